$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 561.2857
$ws.Range("I2").Value = 385.8
$ws.Range("K2").Value = 385.8
$ws.Range("M2").Value = -272.8

$ws.Range("H28").Value = 399.73334
$ws.Range("I28").Value = 368.96295
$ws.Range("J28").Value = 676.6667
$ws.Range("K28").Value = 368.96295
$ws.Range("L28").Value = 676.6667
$ws.Range("M28").Value = 116.03705
$ws.Range("N28").Value = -1646.6667

$ws.Range("H43").Value = 17714
$ws.Range("I43").Value = 99998
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 99998
$ws.Range("L43").Value = 4000
$ws.Range("M43").Value = -99929
$ws.Range("N43").Value = -4138

$ws.Range("H70").Value = 3754
$ws.Range("I70").Value = 3632.6667
$ws.Range("J70").Value = 3794.4443
$ws.Range("K70").Value = 10898.0001
$ws.Range("L70").Value = 11383.3329
$ws.Range("M70").Value = -10628.0001
$ws.Range("N70").Value = -11923.3329

$ws.Range("H73").Value = 3754
$ws.Range("I73").Value = 3632.6667
$ws.Range("J73").Value = 3794.4443
$ws.Range("K73").Value = 10898.0001
$ws.Range("L73").Value = 11383.3329
$ws.Range("M73").Value = -9962.000100000001
$ws.Range("N73").Value = -13255.3329

$ws.Range("H116").Value = 29727.738
$ws.Range("I116").Value = 20915.867
$ws.Range("K116").Value = 20915.867
$ws.Range("M116").Value = -17473.867

$ws.Range("H127").Value = 44761.418
$ws.Range("I127").Value = 61663.293
$ws.Range("J127").Value = 3714
$ws.Range("K127").Value = 184989.879
$ws.Range("L127").Value = 11142
$ws.Range("M127").Value = -180029.879
$ws.Range("N127").Value = -21062

$ws.Range("H138").Value = 3807.849
$ws.Range("I138").Value = 4101.75
$ws.Range("J138").Value = 3360
$ws.Range("K138").Value = 12305.25
$ws.Range("L138").Value = 10080
$ws.Range("M138").Value = -7165.25
$ws.Range("N138").Value = -20360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1149.7307
$ws.Range("I2").Value = 912.2083
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 912.2083
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -799.2083
$ws.Range("N2").Value = -4226

$ws.Range("H61").Value = 2443217.2
$ws.Range("I61").Value = 4279.795
$ws.Range("J61").Value = 50002500
$ws.Range("K61").Value = 4279.795
$ws.Range("L61").Value = 50002500
$ws.Range("M61").Value = -4067.795
$ws.Range("N61").Value = -50002924

$ws.Range("H74").Value = 1211746.4
$ws.Range("I74").Value = 1985671
$ws.Range("J74").Value = 7863.5
$ws.Range("K74").Value = 1985671
$ws.Range("L74").Value = 7863.5
$ws.Range("M74").Value = -1984797
$ws.Range("N74").Value = -9611.5

$ws.Range("H77").Value = 1211746.4
$ws.Range("I77").Value = 1985671
$ws.Range("J77").Value = 7863.5
$ws.Range("K77").Value = 9928355
$ws.Range("L77").Value = 39317.5
$ws.Range("M77").Value = -9923987
$ws.Range("N77").Value = -48053.5

$ws.Range("H102").Value = 4311.375
$ws.Range("I102").Value = 4811.2856
$ws.Range("K102").Value = 4811.2856
$ws.Range("M102").Value = -3189.2856

$ws.Range("H116").Value = 1149.7307
$ws.Range("I116").Value = 912.2083
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 912.2083
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 1381.7917
$ws.Range("N116").Value = -8588

$ws.Range("H132").Value = 864125.75
$ws.Range("I132").Value = 963313.4
$ws.Range("J132").Value = 4499.6665
$ws.Range("K132").Value = 2889940.2
$ws.Range("L132").Value = 13498.9995
$ws.Range("M132").Value = -2887410.2
$ws.Range("N132").Value = -18558.9995

$ws.Range("H136").Value = 2443217.2
$ws.Range("I136").Value = 4279.795
$ws.Range("J136").Value = 50002500
$ws.Range("K136").Value = 12839.385
$ws.Range("L136").Value = 150007500
$ws.Range("M136").Value = -10289.385
$ws.Range("N136").Value = -150012600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1149.7307
$ws.Range("I3").Value = 912.2083
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 912.2083
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -798.2083
$ws.Range("N3").Value = -4228

$ws.Range("H20").Value = 49932.773
$ws.Range("I20").Value = 64252.707
$ws.Range("J20").Value = 1245
$ws.Range("K20").Value = 64252.707
$ws.Range("L20").Value = 1245
$ws.Range("M20").Value = -64005.707
$ws.Range("N20").Value = -1739

$ws.Range("H105").Value = 2279.0557
$ws.Range("I105").Value = 1431.125
$ws.Range("K105").Value = 1431.125
$ws.Range("M105").Value = 315.875

$ws.Range("H107").Value = 1543.6207
$ws.Range("I107").Value = 1491.6072
$ws.Range("K107").Value = 1491.6072
$ws.Range("M107").Value = 428.3928000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2224210.8
$ws.Range("I31").Value = 2316803
$ws.Range("J31").Value = 1995
$ws.Range("K31").Value = 2316803
$ws.Range("L31").Value = 1995
$ws.Range("M31").Value = -2316508
$ws.Range("N31").Value = -2585

$ws.Range("H34").Value = 2224210.8
$ws.Range("I34").Value = 2316803
$ws.Range("J34").Value = 1995
$ws.Range("K34").Value = 2316803
$ws.Range("L34").Value = 1995
$ws.Range("M34").Value = -2316601
$ws.Range("N34").Value = -2399

$ws.Range("H58").Value = 2141572.2
$ws.Range("I58").Value = 2451.04
$ws.Range("J58").Value = 5961431
$ws.Range("K58").Value = 2451.04
$ws.Range("L58").Value = 5961431
$ws.Range("M58").Value = -2248.04
$ws.Range("N58").Value = -5961837

$ws.Range("H132").Value = 1743.3617
$ws.Range("I132").Value = 1542.8096
$ws.Range("J132").Value = 3428
$ws.Range("K132").Value = 4628.4288
$ws.Range("L132").Value = 10284
$ws.Range("M132").Value = -2098.4288
$ws.Range("N132").Value = -15344

$ws.Range("H136").Value = 2141572.2
$ws.Range("I136").Value = 2451.04
$ws.Range("J136").Value = 5961431
$ws.Range("K136").Value = 7353.12
$ws.Range("L136").Value = 17884293
$ws.Range("M136").Value = -4803.12
$ws.Range("N136").Value = -17889393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -20730
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 21000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -20064
$ws.Range("N67").ClearContents()

$ws.Range("H70").Value = 2498.75
$ws.Range("I70").Value = 2331.6667
$ws.Range("K70").Value = 6995.000100000001
$ws.Range("M70").Value = -6680.000100000001

$ws.Range("H73").Value = 2498.75
$ws.Range("I73").Value = 2331.6667
$ws.Range("K73").Value = 6995.000100000001
$ws.Range("M73").Value = -5903.000100000001

$ws.Range("H112").Value = 8166.6665
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws.Range("H113").Value = 1249.8
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 1289.2632
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 3867.7896
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -8207.7896

$ws.Range("H120").Value = 13055.765
$ws.Range("I120").Value = 9149.846
$ws.Range("J120").Value = 25750
$ws.Range("K120").Value = 27449.538
$ws.Range("L120").Value = 77250
$ws.Range("M120").Value = -22611.538
$ws.Range("N120").Value = -86926

$ws.Range("H131").Value = 1894.3864
$ws.Range("I131").Value = 1430
$ws.Range("J131").Value = 1940.825
$ws.Range("K131").Value = 4290
$ws.Range("L131").Value = 5822.475
$ws.Range("M131").Value = 750
$ws.Range("N131").Value = -15902.475

$ws.Range("H139").Value = 4294.9614
$ws.Range("I139").Value = 2980.5293
$ws.Range("K139").Value = 8941.5879
$ws.Range("M139").Value = -3801.5879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2427.3845
$ws.Range("I113").Value = 2427.3845
$ws.Range("K113").Value = 2427.3845
$ws.Range("M113").Value = -257.3845000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 40113172
$ws.Range("I136").Value = 40134800
$ws.Range("K136").Value = 120404400
$ws.Range("M136").Value = -120401850
